$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New task descriptions (column A) ---------------------------------
# NOTE: order matters for shared-string index assignment (matches the
# order the strings were introduced in the source document).

$ws.Range("A20").Value = '[xuly.hpp]
Viết hàm:
- Chuẩn hóa tên tài khoản thành in thường toàn bộ các kí tự, không có khoảng trắng thừa'

$ws.Range("A21").Value = '[khachhang.hpp]
Viết hàm:
- Chuyển cây sang mảng
- Sắp xếp danh sách khách hàng theo mã khách hàng'

$ws.Range("A22").Value = '[khachhang.hpp]
Viết hàm:
- In danh sách khách hàng theo chiều dọc/ngang
- Giải phóng danh sách khách hàng'

$ws.Range("A23").Value = '[khachhang.hpp]
Viết hàm:
- Xóa khách hàng ra khỏi danh sách'

$ws.Range("A19").Value = '[khachhang.hpp]
Viết hàm:
- Nhập và thêm một khách hàng vào danh sách
- Tạo mã khách hàng (một số ngẫu nhiên từ 100 đến 999)
- Kiểm tra mã khách hàng đã tồn tại hay chưa'

$ws.Range("A24").Value = '[khachhang.hpp]
Viết hàm:
- Chỉnh sửa thông tin của một khách hàng trong danh sách khách hàng bằng số điện thoại. Chỉ được thay đổi thông tin tài khoản và mật khẩu'

$ws.Range("A25").Value = '[khachhang.hpp]
Viết hàm:
- Đọc danh sách khách hàng từ file "../File/danhsachkhachhang.txt"'

$ws.Range("A15").Value = '[maytinh.hpp]
Viết hàm:
- Đọc danh sách máy tính từ file "../File/danhsachmaytinh.txt"
- Xóa một máy tính từ danh sách'

# --- Assigned person (column D) ---------------------------------------
$ws.Range("D19").Value = "Xuân Sang"
$ws.Range("D20").Value = "Xuân Lam"
$ws.Range("D21").Value = "Thanh Sang"
$ws.Range("D22").Value = "Hải Sơn"
$ws.Range("D23").Value = "Quốc Thắng"
$ws.Range("D24").Value = "Xuân Sang"
$ws.Range("D25").Value = "Xuân Lam"

# --- Row heights --------------------------------------------------------
$ws.Rows.Item(19).RowHeight = 144
$ws.Rows.Item(21).RowHeight = 83.4
$ws.Rows.Item(22).RowHeight = 79.8
$ws.Rows.Item(23).RowHeight = 63.6
$ws.Rows.Item(24).RowHeight = 103.2
$ws.Rows.Item(25).RowHeight = 81

# --- Selection ------------------------------------------------------------
$ws.Range("E19").Select()
